# Add a "DIC average" summary block (columns M:N) to Sheet1, averaging the
# DIC-value columns (D and J) across the five age groups for each model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (row 3): M3 = "model", N3 = "DIC average" --
$ws.Range("M3").Value = "model"
$ws.Range("N3").Value = "DIC average"

# -- Model rows 4..8 (models 1b..1f): label in M, AVERAGE formula in N --
$models = @("1b", "1c", "1d", "1e", "1f")
for ($i = 0; $i -lt $models.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 13).Value = $models[$i]   # column M = 13
    $ws.Cells.Item($row, 14).Formula = "=AVERAGE(D$row,D$($row+5),D$($row+10),D$($row+15),D$($row+20),J$row,J$($row+5),J$($row+10),J$($row+15),J$($row+20))"   # column N = 14
}

# Row 6 (model "1d") carries the same highlight fill as the rest of that row
# (C6/D6/E6/I6/J6/K6 already use style index 6 -> fillId 7, theme color 9)
$ws.Range("M6").Interior.Color = $ws.Range("C6").Interior.Color
$ws.Range("N6").Interior.Color = $ws.Range("C6").Interior.Color

# -- View tweaks: zoom + active selection --
$excel.ActiveWindow.Zoom = 91
$ws.Range("D3").Select()
